$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the placeholder "Default Message" with the real bot greeting (two lines)
$ws.Range("B2").Value = "Hi my name is SAM. `nHow can I help you today ?"

# Wrap the text so both lines are visible in the cell
$ws.Range("B2").WrapText = $true

# Grow row 2 to fit the two wrapped lines
$ws.Rows.Item(2).RowHeight = 28.8

# Resize the columns to comfortably show the longer header/message text
$ws.Columns.Item(1).ColumnWidth = 13.333333333333334
$ws.Columns.Item(2).ColumnWidth = 48.5

# Move the selection highlight to B3, below the updated message cell
$ws.Range("B3").Select()
